$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Cells.Item(2, 4)
$cD.NumberFormat = "@"
$cD.Value = "27.092.21"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(2, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.87%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(3, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.848.72"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(3, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.39%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(4, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.016"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(4, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.81%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(5, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.014"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(5, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.67%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(6, 4)
$cD.NumberFormat = "@"
$cD.Value = "309.22"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(6, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.05%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(7, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.4763"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(7, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +2.05%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(8, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.3684"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(8, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.87%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(9, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.07237"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(9, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.41%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(10, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.9312"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(10, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.92%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(11, 4)
$cD.NumberFormat = "@"
$cD.Value = "19.85"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(11, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.71%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(12, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.07787"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(12, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.14%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(13, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.849.30"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(13, 5)
$cE.NumberFormat = "@"
$cE.Value = "  -0.14%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(14, 4)
$cD.NumberFormat = "@"
$cD.Value = "5.389"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(14, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +2.05%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(15, 4)
$cD.NumberFormat = "@"
$cD.Value = "6.478"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(15, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.37%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(16, 4)
$cD.NumberFormat = "@"
$cD.Value = "89.26"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(16, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.44%  "
$cE.Style = "Normal"
$cE = $ws.Cells.Item(17, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.85%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(18, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.000008691"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(18, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.26%  "
$cE.Style = "Normal"
$cE = $ws.Cells.Item(19, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.66%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(20, 4)
$cD.NumberFormat = "@"
$cD.Value = "27.100.26"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(20, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.79%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(21, 4)
$cD.NumberFormat = "@"
$cD.Value = "14.60"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(21, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.50%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(22, 4)
$cD.NumberFormat = "@"
$cD.Value = "5.059"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(22, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.81%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(23, 4)
$cD.NumberFormat = "@"
$cD.Value = "10.65"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(23, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.01%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(24, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.939"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(24, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.04%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(25, 4)
$cD.NumberFormat = "@"
$cD.Value = "153.04"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(25, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.41%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(26, 4)
$cD.NumberFormat = "@"
$cD.Value = "18.35"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(26, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.60%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(27, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.986"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(27, 5)
$cE.NumberFormat = "@"
$cE.Value = "  -1.96%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(28, 4)
$cD.NumberFormat = "@"
$cD.Value = "114.72"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(28, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.58%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(29, 4)
$cD.NumberFormat = "@"
$cD.Value = "4.927"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(29, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.93%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(30, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.08872"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(30, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.07%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(31, 4)
$cD.NumberFormat = "@"
$cD.Value = "3.298"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(31, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +2.86%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(32, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.181"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(32, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.07%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(33, 4)
$cD.NumberFormat = "@"
$cD.Value = "4.518"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(33, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.32%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(34, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.7368"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(34, 5)
$cE.NumberFormat = "@"
$cE.Value = "  -1.21%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(35, 4)
$cD.NumberFormat = "@"
$cD.Value = "2.673"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(35, 5)
$cE.NumberFormat = "@"
$cE.Value = "  -3.95%  "
$cE.Style = "Normal"
$cE = $ws.Cells.Item(36, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +3.03%  "
$cE.Style = "Normal"
$cE = $ws.Cells.Item(37, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.78%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(38, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.05266"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(38, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.73%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(39, 4)
$cD.NumberFormat = "@"
$cD.Value = "2.979"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(39, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.42%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(40, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.5277"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(40, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.69%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(41, 4)
$cD.NumberFormat = "@"
$cD.Value = "7.045"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(41, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +2.05%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(42, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.1523"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(42, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.66%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(43, 4)
$cD.NumberFormat = "@"
$cD.Value = "8.292"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(43, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +2.04%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(44, 4)
$cD.NumberFormat = "@"
$cD.Value = "10.61"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(44, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.86%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(45, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.4737"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(45, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.14%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(46, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.014"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(46, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.60%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(47, 4)
$cD.NumberFormat = "@"
$cD.Value = "101.90"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(47, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.51%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(48, 4)
$cD.NumberFormat = "@"
$cD.Value = "1.617"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(48, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.98%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(49, 4)
$cD.NumberFormat = "@"
$cD.Value = "65.77"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(49, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.84%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(50, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.06061"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(50, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +0.47%  "
$cE.Style = "Normal"
$cD = $ws.Cells.Item(51, 4)
$cD.NumberFormat = "@"
$cD.Value = "0.8930"
$cD.Style = "Normal"
$cE = $ws.Cells.Item(51, 5)
$cE.NumberFormat = "@"
$cE.Value = "  +1.27%  "
$cE.Style = "Normal"
